$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.005905376899807375, 0.5206627543158415, 0.486292670446743, 0.6973468795705212, 0.7042605717416123, 51)
    3  = @(0.1582397940106579, 0.7400102157914225, 1.080801010548076, 1.039615799489444, 1.037934161993035, 50)
    4  = @(0.5032026291008911, 0.8922265240678954, 2.069469438159964, 1.438565062192171, 1.361651655305897, 49)
    5  = @(0.2693578958294761, 1.016380599619249, 2.581576023358325, 1.606728360164942, 1.60075163642296, 48)
    6  = @(0.3948146875153234, 0.9953206441242253, 2.563081422080346, 1.600962654805022, 1.568289922150866, 47)
    7  = @(0.2594199654285525, 0.9106554091308247, 2.630675961980921, 1.621935868639978, 1.622546674077204, 38)
    8  = @(0.2821028187081914, 1.011255251399394, 2.812336803878879, 1.677002326736275, 1.675907115328919, 37)
    9  = @(0.1185803502089819, 1.184444939933787, 4.416253151345821, 2.101488318155925, 2.152646315557749, 20)
    10 = @(-0.2292869361971864, 0.9020389666515667, 1.950978095459872, 1.396774174825649, 1.434087187974152, 13)
    11 = @(0.5298242459204994, 0.6893243425020978, 0.5354836367501077, 0.7317674745095656, 0.5643247128038097, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
